# Apply the requested text edits and selection change to the Colouring_tab sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Colouring_tab")

# Update the renamed sample labels (shared strings used in column Q).
$ws.Range("Q2").Value  = "Aci 110"
$ws.Range("Q3").Value  = "A110-1"
$ws.Range("Q7").Value  = "A110-2"
$ws.Range("Q13").Value = "A110-G1"

# Update the active selection on the bottom-right pane to Q2:Q13 (matches the
# saved sheetView selection in the target file).
$ws.Activate()
$ws.Range("Q2:Q13").Select()
